# The workbook lists logging notifications ("Avverkningsanmälningar").
# Column C ("Förändrad" = last-changed date) was bumped for every data
# row (rows 2-348) from Excel serial date 45172 (2023-09-03) to
# 45175 (2023-09-06), i.e. the whole sheet was re-stamped with a new
# "last updated" date during the automatic refresh.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# C2:C348 covers every data row under the header row.
$rng = $ws.Range("C2:C348")
$rng.Value = 45175
